$wb = $excel.ActiveWorkbook

# --- Overview sheet (A1:G7) ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A5").Value = "aa0728e1-c48b-4309-80e6-d9a2b2e0486c.md"
$ws1.Range("B5").Value = "e2e\aa0728e1-c48b-4309-80e6-d9a2b2e0486c.md"
$ws1.Range("E5").Value = "In Translation"
$ws1.Range("F5").Value = "In Translation"
$ws1.Range("G5").Value = "2016-08-29 12:44:21"

$ws1.Range("A6").Value = "a95acbbb-2ea3-4080-844d-5e76f48db359.md"
$ws1.Range("B6").Value = "e2e\a95acbbb-2ea3-4080-844d-5e76f48db359.md"
$ws1.Range("G6").Value = "2016-08-29 12:44:37"

# --- zh-cn sheet (A1:P7) ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A5").Value = "aa0728e1-c48b-4309-80e6-d9a2b2e0486c.md"
$ws2.Range("C5").Value = "In Translation"
$ws2.Range("G5").Value = "aa0728e1-c48b-4309-80e6-d9a2b2e0486c.bc0081baf42b7cd6a64fa875aefbadcfe411bce9.zh-cn.xlf"
$ws2.Range("H5").Value = "2016-08-29 12:44:17"

$ws2.Range("A6").Value = "a95acbbb-2ea3-4080-844d-5e76f48db359.md"
$ws2.Range("G6").Value = "a95acbbb-2ea3-4080-844d-5e76f48db359.421d503524e8baf818e9934c65f19747dc9f0762.zh-cn.xlf"
$ws2.Range("H6").Value = "2016-08-29 12:44:32"

# --- de-de sheet (A1:P7) ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A5").Value = "aa0728e1-c48b-4309-80e6-d9a2b2e0486c.md"
$ws3.Range("C5").Value = "In Translation"
$ws3.Range("G5").Value = "aa0728e1-c48b-4309-80e6-d9a2b2e0486c.bc0081baf42b7cd6a64fa875aefbadcfe411bce9.de-de.xlf"
$ws3.Range("H5").Value = "2016-08-29 12:44:21"

$ws3.Range("A6").Value = "a95acbbb-2ea3-4080-844d-5e76f48db359.md"
$ws3.Range("G6").Value = "a95acbbb-2ea3-4080-844d-5e76f48db359.421d503524e8baf818e9934c65f19747dc9f0762.de-de.xlf"
$ws3.Range("H6").Value = "2016-08-29 12:44:37"
